$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 769713.4399999999
$ws.Range("I92").Value = 879544.1
$ws.Range("J92").Value = 898.5
$ws.Range("K92").Value = 879544.1
$ws.Range("L92").Value = 898.5
$ws.Range("M92").Value = -878296.1
$ws.Range("N92").Value = -3394.5
$ws.Range("H98").Value = 1375.5238
$ws.Range("I98").Value = 1394.3
$ws.Range("K98").Value = 1394.3
$ws.Range("M98").Value = 103.7
$ws.Range("H100").Value = 935.75
$ws.Range("I100").Value = 764.8
$ws.Range("K100").Value = 764.8
$ws.Range("M100").Value = -223.8
$ws.Range("H106").Value = 7669.6665
$ws.Range("I106").Value = 3997
$ws.Range("J106").Value = 9506
$ws.Range("K106").Value = 3997
$ws.Range("L106").Value = 9506
$ws.Range("M106").Value = -3366
$ws.Range("N106").Value = -10768
$ws.Range("H122").Value = 1375.5238
$ws.Range("I122").Value = 1394.3
$ws.Range("K122").Value = 4182.9
$ws.Range("M122").Value = -1732.9
$ws.Range("H129").Value = 890.4533
$ws.Range("I129").Value = 1266.6666
$ws.Range("K129").Value = 3799.9998
$ws.Range("M129").Value = 1200.0002
$ws.Range("H131").Value = 3523.182
$ws.Range("J131").Value = 4385.125
$ws.Range("L131").Value = 13155.375
$ws.Range("N131").Value = -23235.375
$ws.Range("H138").Value = 2735.532
$ws.Range("I138").Value = 2518.2144
$ws.Range("J138").Value = 3055.7896
$ws.Range("K138").Value = 7554.6432
$ws.Range("L138").Value = 9167.3688
$ws.Range("M138").Value = -2414.6432
$ws.Range("N138").Value = -19447.3688

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3535.8655
$ws.Range("I32").Value = 2413.848
$ws.Range("J32").Value = 12138
$ws.Range("K32").Value = 2413.848
$ws.Range("L32").Value = 12138
$ws.Range("M32").Value = -2126.848
$ws.Range("N32").Value = -12712
$ws.Range("H132").Value = 1450.7317
$ws.Range("I132").Value = 969.1739
$ws.Range("K132").Value = 2907.5217
$ws.Range("M132").Value = -377.5217000000002

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 32226.666
$ws.Range("J81").Value = 32226.666
$ws.Range("L81").Value = 32226.666
$ws.Range("N81").Value = -34348.666
$ws.Range("H84").Value = 32226.666
$ws.Range("J84").Value = 32226.666
$ws.Range("L84").Value = 96679.99800000001
$ws.Range("N84").Value = -107287.998
$ws.Range("H86").Value = 78701.96000000001
$ws.Range("I86").Value = 1621.6842
$ws.Range("K86").Value = 1621.6842
$ws.Range("M86").Value = -498.6841999999999
$ws.Range("H89").Value = 78701.96000000001
$ws.Range("I89").Value = 1621.6842
$ws.Range("K89").Value = 8108.420999999999
$ws.Range("M89").Value = -2492.420999999999
$ws.Range("H134").Value = 7715.1797
$ws.Range("I134").Value = 7348.2285
$ws.Range("J134").Value = 10926
$ws.Range("K134").Value = 22044.6855
$ws.Range("L134").Value = 32778
$ws.Range("M134").Value = -19509.6855
$ws.Range("N134").Value = -37848

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 9853.333000000001
$ws.Range("J55").Value = 10280
$ws.Range("L55").Value = 10280
$ws.Range("N55").Value = -10910
$ws.Range("H58").Value = 1554335.6
$ws.Range("I58").Value = 3624242.2
$ws.Range("J58").Value = 1905.6875
$ws.Range("K58").Value = 3624242.2
$ws.Range("L58").Value = 1905.6875
$ws.Range("M58").Value = -3624039.2
$ws.Range("N58").Value = -2311.6875
$ws.Range("H107").Value = 462
$ws.Range("I107").Value = 411.92307
$ws.Range("K107").Value = 411.92307
$ws.Range("M107").Value = 1508.07693
$ws.Range("H132").Value = 2542.3635
$ws.Range("I132").Value = 1536.7778
$ws.Range("K132").Value = 4610.3334
$ws.Range("M132").Value = -2080.3334
$ws.Range("H134").Value = 999.8182
$ws.Range("I134").Value = 999.7778
$ws.Range("K134").Value = 2999.3334
$ws.Range("M134").Value = -464.3334
$ws.Range("H136").Value = 1554335.6
$ws.Range("I136").Value = 3624242.2
$ws.Range("J136").Value = 1905.6875
$ws.Range("K136").Value = 10872726.6
$ws.Range("L136").Value = 5717.0625
$ws.Range("M136").Value = -10870176.6
$ws.Range("N136").Value = -10817.0625

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 820.8889
$ws.Range("J5").Value = 847.3333
$ws.Range("L5").Value = 2541.9999
$ws.Range("N5").Value = -2765.9999
$ws.Range("H14").Value = 147.25
$ws.Range("I14").Value = 147.25
$ws.Range("K14").Value = 441.75
$ws.Range("M14").Value = -268.75
$ws.Range("H103").Value = 2328.4
$ws.Range("I103").Value = 1400.1428
$ws.Range("J103").Value = 4494.3335
$ws.Range("K103").Value = 4200.428400000001
$ws.Range("L103").Value = 13483.0005
$ws.Range("M103").Value = -3321.428400000001
$ws.Range("N103").Value = -15241.0005
$ws.Range("H114").Value = 35717216
$ws.Range("J114").Value = 47622610
$ws.Range("L114").Value = 142867830
$ws.Range("N114").Value = -142874338
$ws.Range("H135").Value = 820.8889
$ws.Range("J135").Value = 847.3333
$ws.Range("L135").Value = 7625.9997
$ws.Range("N135").Value = -12695.9997

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4609
$ws.Range("I70").Value = 4557
$ws.Range("K70").Value = 4557
$ws.Range("M70").Value = -4287
$ws.Range("H73").Value = 4609
$ws.Range("I73").Value = 4557
$ws.Range("K73").Value = 4557
$ws.Range("M73").Value = -3621
$ws.Range("H122").Value = 1762.2142
$ws.Range("I122").Value = 1401.8572
$ws.Range("K122").Value = 4205.571599999999
$ws.Range("M122").Value = -1755.571599999999
$ws.Range("H123").Value = 15479.5
$ws.Range("J123").Value = 15479.5
$ws.Range("L123").Value = 15479.5
$ws.Range("N123").Value = -20379.5
$ws.Range("H135").Value = 27999.4
$ws.Range("J135").Value = 27999.4
$ws.Range("L135").Value = 27999.4
$ws.Range("N135").Value = -38139.4

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3137.8262
$ws.Range("I40").Value = 1093.25
$ws.Range("K40").Value = 1093.25
$ws.Range("M40").Value = -957.25
$ws.Range("H46").Value = 2479.7
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 2621.889
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 2621.889
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -2997.889
$ws.Range("H63").Value = 29110
$ws.Range("J63").Value = 29110
$ws.Range("L63").Value = 29110
$ws.Range("N63").Value = -30608
$ws.Range("H66").Value = 29110
$ws.Range("J66").Value = 29110
$ws.Range("L66").Value = 87330
$ws.Range("N66").Value = -94818
$ws.Range("H93").Value = 459.58334
$ws.Range("I93").Value = 398
$ws.Range("J93").Value = 644.3333
$ws.Range("K93").Value = 398
$ws.Range("L93").Value = 644.3333
$ws.Range("M93").Value = 850
$ws.Range("N93").Value = -3140.3333
$ws.Range("H132").Value = 1974.4642
$ws.Range("I132").Value = 1900.6666
$ws.Range("K132").Value = 5701.9998
$ws.Range("M132").Value = -3171.9998
$ws.Range("H136").Value = 3755.1333
$ws.Range("I136").Value = 1761.1428
$ws.Range("K136").Value = 5283.428400000001
$ws.Range("M136").Value = -2733.428400000001

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1738.258
$ws.Range("I132").Value = 1341.8462
$ws.Range("J132").Value = 3799.6
$ws.Range("K132").Value = 4025.5386
$ws.Range("L132").Value = 11398.8
$ws.Range("M132").Value = -1495.5386
$ws.Range("N132").Value = -16458.8
$ws.Range("H135").Value = 86262.7
$ws.Range("J135").Value = 86262.7
$ws.Range("L135").Value = 86262.7
$ws.Range("N135").Value = -96402.7
